$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(18, 8).Value = 1825.1904
$ws_ALC.Cells.Item(18, 9).Value = 1867.3158
$ws_ALC.Cells.Item(18, 10).Value = 1425
$ws_ALC.Cells.Item(18, 11).Value = 1867.3158
$ws_ALC.Cells.Item(18, 12).Value = 1425
$ws_ALC.Cells.Item(18, 13).Value = -1583.3158
$ws_ALC.Cells.Item(18, 14).Value = -1993
$ws_ALC.Cells.Item(28, 8).Value = 856.5833
$ws_ALC.Cells.Item(28, 9).Value = 567.75
$ws_ALC.Cells.Item(28, 11).Value = 567.75
$ws_ALC.Cells.Item(28, 13).Value = -82.75
$ws_ALC.Cells.Item(74, 8).Value = 6339.3213
$ws_ALC.Cells.Item(74, 10).Value = 7096.067
$ws_ALC.Cells.Item(74, 12).Value = 7096.067
$ws_ALC.Cells.Item(74, 14).Value = -8968.066999999999
$ws_ALC.Cells.Item(77, 8).Value = 6339.3213
$ws_ALC.Cells.Item(77, 10).Value = 7096.067
$ws_ALC.Cells.Item(77, 12).Value = 35480.335
$ws_ALC.Cells.Item(77, 14).Value = -44840.335
$ws_ALC.Cells.Item(116, 8).Value = 11072.818
$ws_ALC.Cells.Item(116, 9).Value = 11200.111
$ws_ALC.Cells.Item(116, 10).Value = 10500
$ws_ALC.Cells.Item(116, 11).Value = 11200.111
$ws_ALC.Cells.Item(116, 12).Value = 10500
$ws_ALC.Cells.Item(116, 13).Value = -7758.111000000001
$ws_ALC.Cells.Item(116, 14).Value = -17384
$ws_ALC.Cells.Item(125, 8).Value = 2620.6
$ws_ALC.Cells.Item(125, 10).Value = 5034.5
$ws_ALC.Cells.Item(125, 12).Value = 45310.5
$ws_ALC.Cells.Item(125, 14).Value = -50230.5
$ws_ALC.Cells.Item(132, 8).Value = 2403.3257
$ws_ALC.Cells.Item(132, 9).Value = 2034.4872
$ws_ALC.Cells.Item(132, 11).Value = 6103.461600000001
$ws_ALC.Cells.Item(132, 13).Value = -3573.461600000001
$ws_ALC.Cells.Item(134, 8).Value = 44949.35
$ws_ALC.Cells.Item(134, 10).Value = 44949.35
$ws_ALC.Cells.Item(134, 12).Value = 44949.35
$ws_ALC.Cells.Item(134, 14).Value = -55089.35
$ws_ALC.Cells.Item(135, 8).Value = 7438.4546
$ws_ALC.Cells.Item(135, 9).Value = 2315
$ws_ALC.Cells.Item(135, 11).Value = 20835
$ws_ALC.Cells.Item(135, 13).Value = -18300
$ws_ALC.Cells.Item(137, 8).Value = 4798.3335
$ws_ALC.Cells.Item(137, 9).Value = 0
$ws_ALC.Cells.Item(137, 10).Value = 4798.3335
$ws_ALC.Cells.Item(137, 11).Value = 0
$ws_ALC.Cells.Item(137, 12).Value = 14395.0005
$ws_ALC.Cells.Item(137, 13).ClearContents()
$ws_ALC.Cells.Item(137, 14).Value = -19495.0005
$ws_ALC.Cells.Item(138, 8).Value = 3123.6375
$ws_ALC.Cells.Item(138, 9).Value = 3021.84
$ws_ALC.Cells.Item(138, 10).Value = 3169.9092
$ws_ALC.Cells.Item(138, 11).Value = 9065.52
$ws_ALC.Cells.Item(138, 12).Value = 9509.7276
$ws_ALC.Cells.Item(138, 13).Value = -3925.52
$ws_ALC.Cells.Item(138, 14).Value = -19789.7276

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(45, 8).Value = 3970.6
$ws_ARM.Cells.Item(45, 9).Value = 3618.5
$ws_ARM.Cells.Item(45, 10).Value = 4498.75
$ws_ARM.Cells.Item(45, 11).Value = 3618.5
$ws_ARM.Cells.Item(45, 12).Value = 4498.75
$ws_ARM.Cells.Item(45, 13).Value = -3241.5
$ws_ARM.Cells.Item(45, 14).Value = -5252.75
$ws_ARM.Cells.Item(74, 8).Value = 1807.5
$ws_ARM.Cells.Item(74, 9).Value = 1313.875
$ws_ARM.Cells.Item(74, 11).Value = 1313.875
$ws_ARM.Cells.Item(74, 13).Value = -439.875
$ws_ARM.Cells.Item(77, 8).Value = 1807.5
$ws_ARM.Cells.Item(77, 9).Value = 1313.875
$ws_ARM.Cells.Item(77, 11).Value = 6569.375
$ws_ARM.Cells.Item(77, 13).Value = -2201.375
$ws_ARM.Cells.Item(82, 8).Value = 67500
$ws_ARM.Cells.Item(82, 9).Value = 35000
$ws_ARM.Cells.Item(82, 11).Value = 35000
$ws_ARM.Cells.Item(82, 13).Value = -34639
$ws_ARM.Cells.Item(85, 8).Value = 67500
$ws_ARM.Cells.Item(85, 9).Value = 35000
$ws_ARM.Cells.Item(85, 11).Value = 35000
$ws_ARM.Cells.Item(85, 13).Value = -33752
$ws_ARM.Cells.Item(94, 8).Value = 0
$ws_ARM.Cells.Item(94, 10).Value = 0
$ws_ARM.Cells.Item(94, 12).Value = 0
$ws_ARM.Cells.Item(94, 14).ClearContents()
$ws_ARM.Cells.Item(121, 8).Value = 0
$ws_ARM.Cells.Item(121, 10).Value = 0
$ws_ARM.Cells.Item(121, 12).Value = 0
$ws_ARM.Cells.Item(121, 14).ClearContents()
$ws_ARM.Cells.Item(134, 8).Value = 100214.5
$ws_ARM.Cells.Item(134, 10).Value = 100214.5
$ws_ARM.Cells.Item(134, 12).Value = 100214.5
$ws_ARM.Cells.Item(134, 14).Value = -110354.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(20, 8).Value = 1799.1333
$ws_BSM.Cells.Item(20, 9).Value = 803.3333
$ws_BSM.Cells.Item(20, 11).Value = 803.3333
$ws_BSM.Cells.Item(20, 13).Value = -556.3333
$ws_BSM.Cells.Item(26, 8).Value = 27746.875
$ws_BSM.Cells.Item(26, 9).Value = 27746.875
$ws_BSM.Cells.Item(26, 11).Value = 27746.875
$ws_BSM.Cells.Item(26, 13).Value = -27454.875
$ws_BSM.Cells.Item(75, 8).Value = 5196.875
$ws_BSM.Cells.Item(75, 9).Value = 5196.875
$ws_BSM.Cells.Item(75, 11).Value = 5196.875
$ws_BSM.Cells.Item(75, 13).Value = -4260.875
$ws_BSM.Cells.Item(78, 8).Value = 5196.875
$ws_BSM.Cells.Item(78, 9).Value = 5196.875
$ws_BSM.Cells.Item(78, 11).Value = 15590.625
$ws_BSM.Cells.Item(78, 13).Value = -10910.625
$ws_BSM.Cells.Item(96, 8).Value = 14499
$ws_BSM.Cells.Item(96, 9).Value = 14499
$ws_BSM.Cells.Item(96, 11).Value = 14499
$ws_BSM.Cells.Item(96, 13).Value = -11753
$ws_BSM.Cells.Item(105, 8).Value = 2387.3333
$ws_BSM.Cells.Item(105, 9).Value = 2422.5454
$ws_BSM.Cells.Item(105, 10).Value = 2000
$ws_BSM.Cells.Item(105, 11).Value = 2422.5454
$ws_BSM.Cells.Item(105, 12).Value = 2000
$ws_BSM.Cells.Item(105, 13).Value = -675.5454
$ws_BSM.Cells.Item(105, 14).Value = -5494
$ws_BSM.Cells.Item(132, 8).Value = 0
$ws_BSM.Cells.Item(132, 10).Value = 0
$ws_BSM.Cells.Item(132, 12).Value = 0
$ws_BSM.Cells.Item(132, 14).ClearContents()

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(22, 8).Value = 453.5
$ws_CRP.Cells.Item(22, 9).Value = 397
$ws_CRP.Cells.Item(22, 11).Value = 397
$ws_CRP.Cells.Item(22, 13).Value = -47
$ws_CRP.Cells.Item(58, 8).Value = 3199.3333
$ws_CRP.Cells.Item(58, 9).Value = 3439.4
$ws_CRP.Cells.Item(58, 10).Value = 1999
$ws_CRP.Cells.Item(58, 11).Value = 3439.4
$ws_CRP.Cells.Item(58, 12).Value = 1999
$ws_CRP.Cells.Item(58, 13).Value = -3236.4
$ws_CRP.Cells.Item(58, 14).Value = -2405
$ws_CRP.Cells.Item(93, 8).Value = 32357.857
$ws_CRP.Cells.Item(93, 10).Value = 69999
$ws_CRP.Cells.Item(93, 12).Value = 69999
$ws_CRP.Cells.Item(93, 14).Value = -73743
$ws_CRP.Cells.Item(122, 8).Value = 2873.5
$ws_CRP.Cells.Item(122, 9).Value = 2873.5
$ws_CRP.Cells.Item(122, 11).Value = 8620.5
$ws_CRP.Cells.Item(122, 13).Value = -6170.5
$ws_CRP.Cells.Item(132, 8).Value = 1634.0667
$ws_CRP.Cells.Item(132, 9).Value = 1741.9166
$ws_CRP.Cells.Item(132, 11).Value = 5225.7498
$ws_CRP.Cells.Item(132, 13).Value = -2695.7498
$ws_CRP.Cells.Item(134, 8).Value = 3332.8125
$ws_CRP.Cells.Item(134, 9).Value = 3332.8125
$ws_CRP.Cells.Item(134, 11).Value = 9998.4375
$ws_CRP.Cells.Item(134, 13).Value = -7463.4375
$ws_CRP.Cells.Item(136, 8).Value = 3199.3333
$ws_CRP.Cells.Item(136, 9).Value = 3439.4
$ws_CRP.Cells.Item(136, 10).Value = 1999
$ws_CRP.Cells.Item(136, 11).Value = 10318.2
$ws_CRP.Cells.Item(136, 12).Value = 5997
$ws_CRP.Cells.Item(136, 13).Value = -7768.200000000001
$ws_CRP.Cells.Item(136, 14).Value = -11097
$ws_CRP.Cells.Item(141, 8).Value = 46666.332
$ws_CRP.Cells.Item(141, 9).Value = 20000
$ws_CRP.Cells.Item(141, 11).Value = 20000
$ws_CRP.Cells.Item(141, 13).Value = -14820

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(63, 8).Value = 0
$ws_CUL.Cells.Item(63, 9).Value = 0
$ws_CUL.Cells.Item(63, 11).Value = 0
$ws_CUL.Cells.Item(63, 13).ClearContents()
$ws_CUL.Cells.Item(66, 8).Value = 0
$ws_CUL.Cells.Item(66, 9).Value = 0
$ws_CUL.Cells.Item(66, 11).Value = 0
$ws_CUL.Cells.Item(66, 13).ClearContents()
$ws_CUL.Cells.Item(70, 8).Value = 9040.5
$ws_CUL.Cells.Item(70, 9).Value = 5387.3335
$ws_CUL.Cells.Item(70, 11).Value = 16162.0005
$ws_CUL.Cells.Item(70, 13).Value = -15847.0005
$ws_CUL.Cells.Item(73, 8).Value = 9040.5
$ws_CUL.Cells.Item(73, 9).Value = 5387.3335
$ws_CUL.Cells.Item(73, 11).Value = 16162.0005
$ws_CUL.Cells.Item(73, 13).Value = -15070.0005
$ws_CUL.Cells.Item(107, 8).Value = 1311.3829
$ws_CUL.Cells.Item(107, 10).Value = 1282.2307
$ws_CUL.Cells.Item(107, 12).Value = 3846.6921
$ws_CUL.Cells.Item(107, 14).Value = -7686.6921
$ws_CUL.Cells.Item(113, 8).Value = 963.9091
$ws_CUL.Cells.Item(113, 9).Value = 1294.25
$ws_CUL.Cells.Item(113, 11).Value = 3882.75
$ws_CUL.Cells.Item(113, 13).Value = -1712.75
$ws_CUL.Cells.Item(134, 8).Value = 9190
$ws_CUL.Cells.Item(134, 9).Value = 1494.5
$ws_CUL.Cells.Item(134, 10).Value = 13999.6875
$ws_CUL.Cells.Item(134, 11).Value = 4483.5
$ws_CUL.Cells.Item(134, 12).Value = 41999.0625
$ws_CUL.Cells.Item(134, 13).Value = 586.5
$ws_CUL.Cells.Item(134, 14).Value = -52139.0625

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(92, 8).Value = 10062.25
$ws_GSM.Cells.Item(92, 10).Value = 10062.25
$ws_GSM.Cells.Item(92, 12).Value = 10062.25
$ws_GSM.Cells.Item(92, 14).Value = -13806.25
$ws_GSM.Cells.Item(126, 8).Value = 3744.1
$ws_GSM.Cells.Item(126, 9).Value = 2492.3333
$ws_GSM.Cells.Item(126, 10).Value = 7499.4
$ws_GSM.Cells.Item(126, 11).Value = 7476.999899999999
$ws_GSM.Cells.Item(126, 12).Value = 22498.2
$ws_GSM.Cells.Item(126, 13).Value = -5006.999899999999
$ws_GSM.Cells.Item(126, 14).Value = -27438.2
$ws_GSM.Cells.Item(135, 8).Value = 83713.14
$ws_GSM.Cells.Item(135, 10).Value = 83713.14
$ws_GSM.Cells.Item(135, 12).Value = 83713.14
$ws_GSM.Cells.Item(135, 14).Value = -93853.14

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(22, 8).Value = 2106.8572
$ws_LTW.Cells.Item(22, 9).Value = 0
$ws_LTW.Cells.Item(22, 11).Value = 0
$ws_LTW.Cells.Item(22, 13).ClearContents()
$ws_LTW.Cells.Item(25, 8).Value = 1950
$ws_LTW.Cells.Item(25, 9).Value = 1950
$ws_LTW.Cells.Item(25, 11).Value = 1950
$ws_LTW.Cells.Item(25, 13).Value = -1720
$ws_LTW.Cells.Item(27, 8).Value = 2106.8572
$ws_LTW.Cells.Item(27, 9).Value = 0
$ws_LTW.Cells.Item(27, 11).Value = 0
$ws_LTW.Cells.Item(27, 13).ClearContents()
$ws_LTW.Cells.Item(48, 8).Value = 43332.668
$ws_LTW.Cells.Item(48, 9).Value = 39999.5
$ws_LTW.Cells.Item(48, 11).Value = 39999.5
$ws_LTW.Cells.Item(48, 13).Value = -39338.5
$ws_LTW.Cells.Item(55, 8).Value = 1368.8667
$ws_LTW.Cells.Item(55, 9).Value = 407.83334
$ws_LTW.Cells.Item(55, 11).Value = 407.83334
$ws_LTW.Cells.Item(55, 13).Value = -234.83334
$ws_LTW.Cells.Item(63, 8).Value = 0
$ws_LTW.Cells.Item(63, 9).Value = 0
$ws_LTW.Cells.Item(63, 11).Value = 0
$ws_LTW.Cells.Item(63, 13).ClearContents()
$ws_LTW.Cells.Item(66, 8).Value = 0
$ws_LTW.Cells.Item(66, 9).Value = 0
$ws_LTW.Cells.Item(66, 11).Value = 0
$ws_LTW.Cells.Item(66, 13).ClearContents()

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(32, 8).Value = 9750
$ws_WVR.Cells.Item(32, 9).Value = 9750
$ws_WVR.Cells.Item(32, 11).Value = 9750
$ws_WVR.Cells.Item(32, 13).Value = -9433
$ws_WVR.Cells.Item(34, 8).Value = 100000
$ws_WVR.Cells.Item(34, 11).Value = 100000
$ws_WVR.Cells.Item(34, 13).Value = -99797
